$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (V1:Y1)
$ws.Range("V1").Value = "ATTA H Cost for Failures"
$ws.Range("W1").Value = "ATTA R Cost for Failures"
$ws.Range("X1").Value = "Tsarouchi MIN H Cost for Failures"
$ws.Range("Y1").Value = "Tsarouchi MIN R Cost for Failures"

# New data values for rows 2-11, columns V-Y
$data = @(
    @(2.6720876217459, 12.3283111690625, 14.4649649718333, 2.8842895980363199),
    @(0.74408058723098702, 18.3688860622863, 20.297275063517802, 3.1949536758237702),
    @(2.6865228560818499, 13.3286065976421, 16.337932925933199, 3.81942119939889),
    @(2.9310042475495801, 13.579284661867799, 20.928532289325801, 2.8118272555038901),
    @(1.7525117617388699, 15.182411123623201, 19.183187749052902, 2.4946626726748899),
    @(1.05604380662411, 16.093658535388499, 15.6912289900819, 3.2342542674729402),
    @(2.55449287289477, 14.920576240820999, 18.670255461831299, 3.7582350153034598),
    @(3.5554878088780999, 14.868615354645099, 21.375834148329599, 2.8382004816587201),
    @(1.3740478201514501, 13.595029701367601, 19.014724432337701, 2.2979210724330299),
    @(2.04408086108378, 15.5826476926888, 16.418468265158801, 3.8498519586645799)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("V$row").Value = $data[$i][0]
    $ws.Range("W$row").Value = $data[$i][1]
    $ws.Range("X$row").Value = $data[$i][2]
    $ws.Range("Y$row").Value = $data[$i][3]
}

# Update sheet view selection to match the post-edit state
$ws.Range("L3").Select() | Out-Null
